# Integrate Agora data for hydgn/BHPSbP
# Replace the "U.S. DOE Fuel Cell Technologies Office" source with the
# "European Commission" source (Hydrogen generation in Europe report).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("BHPSbP")

# --- Update the "About" sheet source citation (rows 3-7) ---
$ws1.Range("B3").Value = "European Commission"
$ws1.Range("B4").Value = 2020
$ws1.Range("B6").Value = "https://op.europa.eu/en/publication-detail/-/publication/7e4afa7d-d077-11ea-adf7-01aa75ed71a1"
$ws1.Range("B5").Value = "Hydrogen generation in Europe: Overview of costs and key benefits"
$ws1.Range("B7").Value = "Introduction (paragraphs 3 and 4)"

# --- Restore view selections to match the saved workbook state ---
$ws2.Range("B6").Select()
$ws1.Range("B8").Select()
